# space-type: add acronym column in template
#
# The template's table currently reads:  name | description
# It becomes:                            name | acronym | description
#
# i.e. a new "acronym" column is inserted as column B, and the existing
# "description" column shifts from B to C (keeping its own data/values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-home the existing "description" column into the new column C,
#    carrying over its header text and all of its data values.
$ws.Range("C1").Value = "description"
$ws.Range("C2").Value = "Good description"
$ws.Range("C3").Value = "Bad description"
$ws.Range("C4").Value = "Medium description"

# The "description" header keeps the same bold/size-11 look that the
# "name" header already uses.
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Font.Size = 11

# 2) Put the new "acronym" column into B. Its header keeps the bold
#    size-12 styling that used to belong to the "description" header.
$ws.Range("B1").Value = "acronym"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 12

$ws.Range("B2").Value = "o"
$ws.Range("B3").Value = "t"
$ws.Range("B4").Value = "t"

# 3) Size the new description column (C) to fit its text, mirroring the
#    original author's "AutoFit Column Width".
$ws.Columns.Item(3).AutoFit()

# 4) Set the page to portrait orientation (as in the final template) and
#    leave the selection where the author left off editing the table.
$ws.PageSetup.Orientation = 1
$ws.Range("B6").Select()
